# Applies the diff:
#  1) Column C (the "Förändrad" / changed date) goes from 45207 to 45208
#     for every data row (rows 2..173).
#  2) For rows 2..5, the hyperlink formulas in columns S, T, V, W, X, Y
#     have "Logging_ALVDALEN" replaced with "Logging_2039" in their URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 173

# 1) Update the "Förändrad" date column (C) for all data rows.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}

# 2) Rename the "Logging_ALVDALEN" folder segment to "Logging_2039"
#    inside the HYPERLINK formulas for rows 2-5, columns S,T,V,W,X,Y.
$hlCols = @("S", "T", "V", "W", "X", "Y")
for ($r = 2; $r -le 5; $r++) {
    foreach ($col in $hlCols) {
        $rng = $ws.Range($col + $r)
        $f = $rng.Formula
        if ($f -ne $null -and $f -like "*Logging_ALVDALEN*") {
            $rng.Formula = $f.Replace("Logging_ALVDALEN", "Logging_2039")
        }
    }
}
